$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1898.5
$ws.Range("I28").Value = 271.33334
$ws.Range("J28").Value = 3525.6667
$ws.Range("K28").Value = 271.33334
$ws.Range("L28").Value = 3525.6667
$ws.Range("M28").Value = 213.66666
$ws.Range("N28").Value = -4495.6667
$ws.Range("H69").Value = 1809.3334
$ws.Range("I69").Value = 1456.5
$ws.Range("K69").Value = 4369.5
$ws.Range("M69").Value = -3495.5
$ws.Range("H72").Value = 1809.3334
$ws.Range("I72").Value = 1456.5
$ws.Range("K72").Value = 13108.5
$ws.Range("M72").Value = -8740.5
$ws.Range("H137").Value = 942.2
$ws.Range("I137").Value = 942.2
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2826.6
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -276.6000000000004
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 6217.645
$ws.Range("J138").Value = 6618.6665
$ws.Range("L138").Value = 19855.9995
$ws.Range("N138").Value = -30135.9995
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 200620
$ws.Range("I14").Value = 200620
$ws.Range("K14").Value = 200620
$ws.Range("M14").Value = -200445
$ws.Range("H32").Value = 10427.546
$ws.Range("I32").Value = 10170.3
$ws.Range("K32").Value = 10170.3
$ws.Range("M32").Value = -9883.299999999999
$ws.Range("H45").Value = 2499.8
$ws.Range("I45").Value = 2624.75
$ws.Range("K45").Value = 2624.75
$ws.Range("M45").Value = -2247.75
$ws.Range("H74").Value = 2212.6667
$ws.Range("I74").Value = 1812
$ws.Range("J74").Value = 3014
$ws.Range("K74").Value = 1812
$ws.Range("L74").Value = 3014
$ws.Range("M74").Value = -938
$ws.Range("N74").Value = -4762
$ws.Range("H77").Value = 2212.6667
$ws.Range("I77").Value = 1812
$ws.Range("J77").Value = 3014
$ws.Range("K77").Value = 9060
$ws.Range("L77").Value = 15070
$ws.Range("M77").Value = -4692
$ws.Range("N77").Value = -23806
$ws.Range("H102").Value = 52501428
$ws.Range("I102").Value = 3335236.8
$ws.Range("K102").Value = 3335236.8
$ws.Range("M102").Value = -3333614.8
$ws.Range("H122").Value = 5857.2383
$ws.Range("I122").Value = 2916.8333
$ws.Range("J122").Value = 9777.777
$ws.Range("K122").Value = 8750.499899999999
$ws.Range("L122").Value = 29333.331
$ws.Range("M122").Value = -6300.499899999999
$ws.Range("N122").Value = -34233.331
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1000
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H107").Value = 39071.637
$ws.Range("I107").Value = 51598.875
$ws.Range("K107").Value = 51598.875
$ws.Range("M107").Value = -49678.875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9333.166999999999
$ws.Range("J31").Value = 9333.166999999999
$ws.Range("L31").Value = 9333.166999999999
$ws.Range("N31").Value = -9923.166999999999
$ws.Range("H34").Value = 9333.166999999999
$ws.Range("J34").Value = 9333.166999999999
$ws.Range("L34").Value = 9333.166999999999
$ws.Range("N34").Value = -9737.166999999999
$ws.Range("H51").Value = 32500
$ws.Range("I51").Value = 20000
$ws.Range("J51").Value = 45000
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 45000
$ws.Range("M51").Value = -19264
$ws.Range("N51").Value = -46472
$ws.Range("H61").Value = 32500
$ws.Range("I61").Value = 20000
$ws.Range("J61").Value = 45000
$ws.Range("K61").Value = 20000
$ws.Range("L61").Value = 45000
$ws.Range("M61").Value = -19652
$ws.Range("N61").Value = -45696
$ws.Range("H99").Value = 3043271.5
$ws.Range("I99").Value = 2632962.8
$ws.Range("J99").Value = 3125333.2
$ws.Range("K99").Value = 2632962.8
$ws.Range("L99").Value = 3125333.2
$ws.Range("M99").Value = -2631464.8
$ws.Range("N99").Value = -3128329.2
$ws.Range("H126").Value = 3043271.5
$ws.Range("I126").Value = 2632962.8
$ws.Range("J126").Value = 3125333.2
$ws.Range("K126").Value = 7898888.399999999
$ws.Range("L126").Value = 9375999.600000001
$ws.Range("M126").Value = -7896418.399999999
$ws.Range("N126").Value = -9380939.600000001
$ws.Range("H141").Value = 489387.16
$ws.Range("J141").Value = 489387.16
$ws.Range("L141").Value = 489387.16
$ws.Range("N141").Value = -499747.16
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2078.3333
$ws.Range("I80").Value = 1903.5
$ws.Range("K80").Value = 5710.5
$ws.Range("M80").Value = -4774.5
$ws.Range("H83").Value = 2078.3333
$ws.Range("I83").Value = 1903.5
$ws.Range("K83").Value = 17131.5
$ws.Range("M83").Value = -12451.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H75").Value = 51363.637
$ws.Range("J75").Value = 51363.637
$ws.Range("L75").Value = 51363.637
$ws.Range("N75").Value = -53111.637
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H78").Value = 51363.637
$ws.Range("J78").Value = 51363.637
$ws.Range("L78").Value = 154090.911
$ws.Range("N78").Value = -162826.911
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -830
$ws.Range("H35").Value = 10598.077
$ws.Range("I35").Value = 3968.5715
$ws.Range("K35").Value = 3968.5715
$ws.Range("M35").Value = -3632.5715
$ws.Range("H100").Value = 2970.3
$ws.Range("J100").Value = 1800
$ws.Range("L100").Value = 1800
$ws.Range("N100").Value = -2882
$ws.Range("H132").Value = 3590.3635
$ws.Range("I132").Value = 3590.3635
$ws.Range("K132").Value = 10771.0905
$ws.Range("M132").Value = -8241.0905
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 20000
$ws.Range("I24").Value = 20000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -19770
